# MU.xlsx TanzyWatch update — "10Th - MB for single stock and added new group"
#
# Inserts 3 new columns (C,D,E) right after the existing "B" column, shifting the
# old C/D/E (the Jun_15 / Jun_13 / Jun_10 data) out to F/G/H. The new block gets a
# fresh pair of header dates (Jun_27 in B, Jun_26 spanning C&D) and is seeded with
# the same "UN" placeholder used everywhere else, except for a handful of analysts
# who published a new note on 6/21 (and a couple on 6/18 and 6/22) — those get the
# note text written into both C and D, with the highlighted fill style used
# elsewhere in the sheet for "new" notes. Finally two new analyst rows are appended
# for brand-new coverage (Benchmark, Evercore ISI).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room: insert 3 new columns after B. This shifts the old C/D/E (and
#    their values/styles) to F/G/H automatically.
# ---------------------------------------------------------------------------
$ws.Columns("C:E").Insert()

# ---------------------------------------------------------------------------
# 2. Header row — new dates for the newly inserted columns.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# ---------------------------------------------------------------------------
# 3. Seed every data row (2-27) with the "UN" placeholder in the 3 new columns,
#    matching the rest of the sheet.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"   # C
    $ws.Cells.Item($r, 4).Value = "UN"   # D
    $ws.Cells.Item($r, 5).Value = "UN"   # E
}

# ---------------------------------------------------------------------------
# 4. Copy the sheet's existing "highlighted note" format (from the old E6,
#    now shifted to H6) so new highlighted notes reuse the same fill style
#    instead of inventing a new one.
# ---------------------------------------------------------------------------
$ws.Range("H6").Copy()

# New analyst notes dated 6/21/2018 — written into both C & D, highlighted.
$newNotes = @{
    3  = "6/21/2018,Set Price Target,Outperform -> Buy,`$80.00 -> `$90.00"
    14 = "6/21/2018,Raises Target,Buy -> Buy,`$108.00"
    17 = "6/21/2018,Raises Target,Overweight -> Buy,`$82.00 -> `$84.00"
    18 = "6/21/2018,Raises Target,Outperform -> Positive,`$83.00"
    27 = "6/21/2018,Raises Target,Overweight,`$67.00"
}

foreach ($r in $newNotes.Keys) {
    $ws.Range("C$r").PasteSpecial(-4122)
    $ws.Range("D$r").PasteSpecial(-4122)
    $ws.Range("C$r").Value = $newNotes[$r]
    $ws.Range("D$r").Value = $newNotes[$r]
}

# ---------------------------------------------------------------------------
# 5. Two brand-new analysts covering the stock, appended as new rows.
# ---------------------------------------------------------------------------

# Row 28: Benchmark — Initiates coverage (no highlight, default style).
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "6/22/2018,Initiates,Buy,`$80.00"
$ws.Range("D28").Value = "6/22/2018,Initiates,Buy,`$80.00"

# Row 29: Evercore ISI — Raises Target (highlighted, like the other new notes).
$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("C29").Value = "6/18/2018,Raises Target,Outperform,`$80.00 -> `$100.00"
$ws.Range("D29").Value = "6/18/2018,Raises Target,Outperform,`$80.00 -> `$100.00"

# ---------------------------------------------------------------------------
# 6. Match column widths for the newly inserted columns to the rest of the
#    narrow "data" columns in the block.
# ---------------------------------------------------------------------------
$ws.Range("C1:E1").ColumnWidth = 8
